$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows before the current row 12 ("Programa resumido:"),
# shifting everything from row 12 onward down to row 15 onward.
$ws.Rows("12:14").Insert()

# Row 12: new section header "Docentes responsáveis:" (column A only)
$ws.Range("A12").Value = "Docentes responsáveis:"

# Row 13: first responsible professor (columns B and C)
$ws.Range("B13").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C13").Value = "5840726 - Cristina Bormio Nunes"

# Row 14: second responsible professor (columns B and C)
$ws.Range("B14").Value = "1341653 - Maria José Ramos Sandim"
$ws.Range("C14").Value = "1341653 - Maria José Ramos Sandim"

# Remove the leftover empty cells that Insert() duplicated from the row above,
# so the sparse cell layout matches rows such as "Avaliação:" / data-only rows.
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
